$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update EMPRESA column (E) from "GUANABARA" to "ITAPEMIRIM" for the
# specific SOBE/QUI, SOBE/SEX, SOBE/SEG and SOBE/TER rows.
$rows = @(24, 25, 28, 29, 32, 33, 36, 37)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "ITAPEMIRIM"
}

# Scroll / selection state: freeze pane top-left moved to A17 and the
# active selection moved to E38.
$ws.Activate()
$ws.Range("E38").Select()
$excel.ActiveWindow.ScrollRow = 17
